# Scheduled runner update: refresh market-price derived figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4012002.8
$ws.Range("J17").Value = 4012002.8
$ws.Range("L17").Value = 12036008.4
$ws.Range("N17").Value = -12036344.4
$ws.Range("H34").Value = 7531.3
$ws.Range("I34").Value = 3152.1667
$ws.Range("K34").Value = 3152.1667
$ws.Range("M34").Value = -2949.1667
$ws.Range("H36").Value = 7531.3
$ws.Range("I36").Value = 3152.1667
$ws.Range("K36").Value = 3152.1667
$ws.Range("M36").Value = -2437.1667
$ws.Range("H80").Value = 6332.2354
$ws.Range("I80").Value = 468.75
$ws.Range("J80").Value = 20404.6
$ws.Range("K80").Value = 1406.25
$ws.Range("L80").Value = 61213.8
$ws.Range("M80").Value = -408.25
$ws.Range("N80").Value = -63209.8
$ws.Range("H83").Value = 6332.2354
$ws.Range("I83").Value = 468.75
$ws.Range("J83").Value = 20404.6
$ws.Range("K83").Value = 4218.75
$ws.Range("L83").Value = 183641.4
$ws.Range("M83").Value = 773.25
$ws.Range("N83").Value = -193625.4
$ws.Range("H132").Value = 3664491.5
$ws.Range("I132").Value = 1160.4359
$ws.Range("J132").Value = 25644478
$ws.Range("K132").Value = 3481.3077
$ws.Range("L132").Value = 76933434
$ws.Range("M132").Value = -951.3076999999998
$ws.Range("N132").Value = -76938494
$ws.Range("H137").Value = 1323.0303
$ws.Range("I137").Value = 1187.75
$ws.Range("J137").Value = 2080.6
$ws.Range("K137").Value = 3563.25
$ws.Range("L137").Value = 6241.799999999999
$ws.Range("M137").Value = -1013.25
$ws.Range("N137").Value = -11341.8
$ws.Range("H138").Value = 2970.5876
$ws.Range("I138").Value = 1470.8
$ws.Range("J138").Value = 3817.242
$ws.Range("K138").Value = 4412.4
$ws.Range("L138").Value = 11451.726
$ws.Range("M138").Value = 727.6000000000004
$ws.Range("N138").Value = -21731.726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4155
$ws.Range("I15").Value = 1100
$ws.Range("K15").Value = 1100
$ws.Range("M15").Value = -750
$ws.Range("H32").Value = 5851.03
$ws.Range("I32").Value = 5183.1797
$ws.Range("J32").Value = 11254.546
$ws.Range("K32").Value = 5183.1797
$ws.Range("L32").Value = 11254.546
$ws.Range("M32").Value = -4896.1797
$ws.Range("N32").Value = -11828.546
$ws.Range("H61").Value = 503898.84
$ws.Range("I61").Value = 3912
$ws.Range("J61").Value = 1670534.9
$ws.Range("K61").Value = 3912
$ws.Range("L61").Value = 1670534.9
$ws.Range("M61").Value = -3700
$ws.Range("N61").Value = -1670958.9
$ws.Range("H74").Value = 9805290
$ws.Range("I74").Value = 1082.5
$ws.Range("J74").Value = 20001666
$ws.Range("K74").Value = 1082.5
$ws.Range("L74").Value = 20001666
$ws.Range("M74").Value = -208.5
$ws.Range("N74").Value = -20003414
$ws.Range("H77").Value = 9805290
$ws.Range("I77").Value = 1082.5
$ws.Range("J77").Value = 20001666
$ws.Range("K77").Value = 5412.5
$ws.Range("L77").Value = 100008330
$ws.Range("M77").Value = -1044.5
$ws.Range("N77").Value = -100017066
$ws.Range("H132").Value = 2636604.2
$ws.Range("I132").Value = 2853.5652
$ws.Range("J132").Value = 6675022
$ws.Range("K132").Value = 8560.695599999999
$ws.Range("L132").Value = 20025066
$ws.Range("M132").Value = -6030.695599999999
$ws.Range("N132").Value = -20030126
$ws.Range("H136").Value = 503898.84
$ws.Range("I136").Value = 3912
$ws.Range("J136").Value = 1670534.9
$ws.Range("K136").Value = 11736
$ws.Range("L136").Value = 5011604.699999999
$ws.Range("M136").Value = -9186
$ws.Range("N136").Value = -5016704.699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 49145.637
$ws.Range("I134").Value = 3953
$ws.Range("J134").Value = 202800.6
$ws.Range("K134").Value = 11859
$ws.Range("L134").Value = 608401.8
$ws.Range("M134").Value = -9324
$ws.Range("N134").Value = -613471.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4026.9333
$ws.Range("I94").Value = 3181.5
$ws.Range("J94").Value = 4334.364
$ws.Range("K94").Value = 3181.5
$ws.Range("L94").Value = 4334.364
$ws.Range("M94").Value = -2730.5
$ws.Range("N94").Value = -5236.364
$ws.Range("H134").Value = 13208544
$ws.Range("J134").Value = 3334041.8
$ws.Range("L134").Value = 10002125.4
$ws.Range("N134").Value = -10007195.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4722.241
$ws.Range("I122").Value = 996.2353000000001
$ws.Range("J122").Value = 10000.75
$ws.Range("K122").Value = 8966.117700000001
$ws.Range("L122").Value = 90006.75
$ws.Range("M122").Value = -6516.117700000001
$ws.Range("N122").Value = -94906.75
$ws.Range("H131").Value = 2042008.2
$ws.Range("I131").Value = 5882980.5
$ws.Range("J131").Value = 1491.7812
$ws.Range("K131").Value = 17648941.5
$ws.Range("L131").Value = 4475.3436
$ws.Range("M131").Value = -17643901.5
$ws.Range("N131").Value = -14555.3436
$ws.Range("H132").Value = 4631979.5
$ws.Range("I132").Value = 1369.3
$ws.Range("J132").Value = 7939558.5
$ws.Range("K132").Value = 12323.7
$ws.Range("L132").Value = 71456026.5
$ws.Range("M132").Value = -9793.699999999999
$ws.Range("N132").Value = -71461086.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11117.857
$ws.Range("I126").Value = 14232.667
$ws.Range("J126").Value = 3330.8333
$ws.Range("K126").Value = 42698.001
$ws.Range("L126").Value = 9992.499899999999
$ws.Range("M126").Value = -40228.001
$ws.Range("N126").Value = -14932.4999
$ws.Range("H132").Value = 12826482
$ws.Range("I132").Value = 18520364
$ws.Range("J132").Value = 15249.75
$ws.Range("K132").Value = 55561092
$ws.Range("L132").Value = 45749.25
$ws.Range("M132").Value = -55558562
$ws.Range("N132").Value = -50809.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6791855
$ws.Range("I122").Value = 8936889
$ws.Range("J122").Value = 2501787.5
$ws.Range("K122").Value = 26810667
$ws.Range("L122").Value = 7505362.5
$ws.Range("M122").Value = -26808217
$ws.Range("N122").Value = -7510262.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1389.6666
$ws.Range("I122").Value = 1358.1428
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4074.4284
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1624.4284
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 1423.4131
$ws.Range("I132").Value = 552.9375
$ws.Range("J132").Value = 3413.0715
$ws.Range("K132").Value = 1658.8125
$ws.Range("L132").Value = 10239.2145
$ws.Range("M132").Value = 871.1875
$ws.Range("N132").Value = -15299.2145
